# Update template files for import trips
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text changes: mark the two columns as required with a trailing "*"
$ws.Range("F1").Value = "Needs Delivery Note ?*"
$ws.Range("G1").Value = "Has Attchment ?*"

# Row 2: replace the sample pickup-window dates with plain text dates and
# clear the "end date" column
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "10/21/2021"
$ws.Range("C2").ClearContents() | Out-Null

# Row 3: same treatment
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "10/22/2021"
$ws.Range("C3").ClearContents() | Out-Null

# Constrain the Yes/No columns to a dropdown list (xlValidateList, xlValidAlertStop)
$ws.Range("F2:G3").Validation.Add(3, 1, 1, """yes,no""") | Out-Null

# Move the active selection from H3 to F3
$ws.Range("F3").Select() | Out-Null
